$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 413.13
$ws.Range("L2").Value = -989969.67
$ws.Range("N2").Value = -19726.16
$ws.Range("P2").Value = -1061202
$ws.Range("R2").Value = -3083621.63
$ws.Range("B3").Value = 98.47
$ws.Range("L3").Value = 287017.77
$ws.Range("N3").Value = 15209.64
$ws.Range("P3").Value = 284915.98
$ws.Range("R3").Value = 413128.48
$ws.Range("B4").Value = 43.269
$ws.Range("L4").Value = -702951.9
$ws.Range("N4").Value = -4516.53
$ws.Range("P4").Value = -776286.02
$ws.Range("R4").Value = -2670493.16
$ws.Range("B5").Value = 65.117
$ws.Range("L7").Value = 165
$ws.Range("B8").Value = 35.34
$ws.Range("D8").Value = 17.303
$ws.Range("F8").Value = 65.12
$ws.Range("L8").Value = 101
$ws.Range("F9").Value = 65.12
$ws.Range("L9").Value = 161
$ws.Range("D10").Value = 0
$ws.Range("F10").Value = 65.12
$ws.Range("L10").Value = 106
$ws.Range("B11").Value = 98.47
$ws.Range("F11").Value = 65.12
$ws.Range("L11").Value = 154
$ws.Range("B14").Value = 287017.77
$ws.Range("D14").Value = 15209.64
$ws.Range("F14").Value = 284915.98
$ws.Range("H14").Value = 413128.48
$ws.Range("L14").Value = 1146.33
$ws.Range("M14").Value = 423.93
$ws.Range("N14").Value = 1108.05
$ws.Range("O14").Value = 429.22
$ws.Range("B15").Value = 126110.71
$ws.Range("D15").Value = 397918.84
$ws.Range("F15").Value = 128212.5
$ws.Range("L15").Value = 7944.79
$ws.Range("M15").Value = 30802.14
$ws.Range("N15").Value = 8400.25
$ws.Range("P15").Value = 79142.33
$ws.Range("F16").Value = 21684.02
$ws.Range("H16").Value = 449447.89
$ws.Range("M16").Value = 122049.98
$ws.Range("P16").Value = 340042.18
$ws.Range("B17").Value = 537325.96
$ws.Range("D17").Value = 29316.95
$ws.Range("H17").Value = 1614829.22
$ws.Range("M17").Value = 8610.530000000001
$ws.Range("N17").Value = 456.29
$ws.Range("P17").Value = 25877.29
$ws.Range("B18").Value = 17.3
$ws.Range("D18").Value = 2686.09
$ws.Range("F18").Value = 0
$ws.Range("M18").Value = 362.68
$ws.Range("O18").Value = 0
$ws.Range("D19").Value = 14788.6
$ws.Range("O20").Value = 5057.65
$ws.Range("P20").Value = 39604.92
$ws.Range("M21").Value = -3.98
$ws.Range("N21").Value = -8.800000000000001
$ws.Range("O21").Value = -4.33
$ws.Range("P21").Value = -6.04
$ws.Range("B22").Value = 69.47
$ws.Range("D22").Value = 3.68
$ws.Range("F22").Value = 68.97
$ws.Range("M22").Value = -2.49
$ws.Range("N22").Value = -5.5
$ws.Range("O22").Value = -2.71
$ws.Range("P22").Value = -3.77
$ws.Range("B23").Value = 30.53
$ws.Range("D23").Value = 96.31999999999999
$ws.Range("F23").Value = 31.03
$ws.Range("B24").Value = 1241.83
$ws.Range("D24").Value = 67.76000000000001
$ws.Range("F24").Value = 1326.55
$ws.Range("H24").Value = 3732.08
$ws.Range("B25").Value = 0.04
$ws.Range("D25").Value = 6.21
$ws.Range("F25").Value = 0
$ws.Range("D26").Value = 34.18
